# Generate Report for Handoff
# The "8b1fa46e-9444-4d3b-a2ec-76379a71a016" file moves from "Handed back: in
# sync with en-US" (top of its sheet) to "Ready for handoff" (moved after the
# other two in-sync files, right before the .localization-config row), with a
# freshly generated handoff timestamp.  The hyperlink target URLs are left
# exactly as they were (only their on-screen text / the cell text changes).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "ffff07db2363-c39d-4576-94cf-50a71b3b573b.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"

$ws1.Range("A4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = ".localization-config"
$ws1.Range("B5").Value = "Not to be localized"
$ws1.Range("C5").Value = "Not to be localized"

$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/8b1fa46e-9444-4d3b-a2ec-76379a71a016.md", "", "", "ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md", "", "", "ffff07db2363-c39d-4576-94cf-50a71b3b573b.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/ffff07db2363-c39d-4576-94cf-50a71b3b573b.md", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-03 13:05:40"
$ws2.Range("E2").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md"
$ws2.Range("F2").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-03 13:06:42"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "ffff07db2363-c39d-4576-94cf-50a71b3b573b.md"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-03 13:05:40"
$ws2.Range("E3").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md"
$ws2.Range("F3").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-03-03 13:06:42"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-03 13:11:08"
$ws2.Range("E4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md"
$ws2.Range("F4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-03-03 13:09:49"
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = ".localization-config"
$ws2.Range("B5").Value = "Not to be localized"
$ws2.Range("D5").Value = "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Ignored"

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/8b1fa46e-9444-4d3b-a2ec-76379a71a016.md", "", "", "ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6dfa388f8bad181236f8c32283901eb4690eeca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.zh-cn.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/cfb91711088aaabe3ca773f218cc4e86897315ff/e2e/8b1fa46e-9444-4d3b-a2ec-76379a71a016.md", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7730832365fc9e8d86f8f3073013475bc152e34c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.zh-cn.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md", "", "", "ffff07db2363-c39d-4576-94cf-50a71b3b573b.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4452c0e6c0181e9871e5e56857058a17a95fb3c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fb48695bd52ea3d0029cbbe998ccc534e9d0b022/e2e/f9a1c574-bb71-4c1c-9ebc-c71355efee93.md", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/669b4126d8bedb40bf00e89d5257e7eea5691f92/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/ffff07db2363-c39d-4576-94cf-50a71b3b573b.md", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4452c0e6c0181e9871e5e56857058a17a95fb3c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fb48695bd52ea3d0029cbbe998ccc534e9d0b022/e2e/f9a1c574-bb71-4c1c-9ebc-c71355efee93.md", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/669b4126d8bedb40bf00e89d5257e7eea5691f92/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.zh-cn.xlf", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-03 13:05:52"
$ws3.Range("E2").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md"
$ws3.Range("F2").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-03 13:07:08"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "ffff07db2363-c39d-4576-94cf-50a71b3b573b.md"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-03 13:05:52"
$ws3.Range("E3").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md"
$ws3.Range("F3").Value = "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf"
$ws3.Range("G3").Value = "2016-03-03 13:07:08"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-03 13:11:21"
$ws3.Range("E4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md"
$ws3.Range("F4").Value = "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.de-de.xlf"
$ws3.Range("G4").Value = "2016-03-03 13:10:15"
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = ".localization-config"
$ws3.Range("B5").Value = "Not to be localized"
$ws3.Range("D5").Value = "0001-01-01 00:00:00"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Ignored"

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/8b1fa46e-9444-4d3b-a2ec-76379a71a016.md", "", "", "ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89de05c8dfd99024c461a6d7d79108c307480f77/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.de-de.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/40f6570d504ad74c58412b18dd88515c10411d62/e2e/8b1fa46e-9444-4d3b-a2ec-76379a71a016.md", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a90e122bd3154d524648b2abe63fbb325941959b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.de-de.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/ffffff9c68a8cd-7eeb-4ce7-9284-457673a081d3.md", "", "", "ffff07db2363-c39d-4576-94cf-50a71b3b573b.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/099229d874301cc9884b47531f87901107a53dfd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5e81866c70dcc817b975b0c11de5640cbc0ef19a/e2e/f9a1c574-bb71-4c1c-9ebc-c71355efee93.md", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.md")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/560d4e895349441f4e6bf3aa830e30dc149f3f85/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf", "", "", "f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/e2e/ffff07db2363-c39d-4576-94cf-50a71b3b573b.md", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/099229d874301cc9884b47531f87901107a53dfd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5e81866c70dcc817b975b0c11de5640cbc0ef19a/e2e/f9a1c574-bb71-4c1c-9ebc-c71355efee93.md", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.md")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/560d4e895349441f4e6bf3aa830e30dc149f3f85/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f9a1c574-bb71-4c1c-9ebc-c71355efee93.4ff5f8090bda51aa16bf792a4318619c390bfc28.de-de.xlf", "", "", "8b1fa46e-9444-4d3b-a2ec-76379a71a016.2aeb67eb68ddb87c39197543e27d05f0a24978d5.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2e9d665b8e169570fba59a5d6dd13089f67ff626/.localization-config", "", "", ".localization-config")
